# Bug List documentation workbook — v1.0.4 update
# - Order Status widget shipped along the top of the display (v1.0.4)
# - Record the new version / git hash on the Versions sheet
# - Record which version each Enhancement request shipped in
# - "Bug List" becomes the active/selected tab again

$wb = $excel.ActiveWorkbook

$wsBugs = $wb.Worksheets.Item("Bug List")
$wsEnh  = $wb.Worksheets.Item("Enhancements")
$wsVer  = $wb.Worksheets.Item("Versions")

# --- Versions sheet: git hash for 1.0.3, then the new 1.0.4 row ------------
$wsVer.Range("B5").Value = "6410c4d660bd46c7dd1254e7d8204dd9235d92ef"

# --- Enhancements sheet: note which release each item shipped in -----------
$wsEnh.Range("C2").Value = "1.0.3"
$wsEnh.Range("C3").Value = "1.0.4"
[void]$wsEnh.Range("C3").Select()

$wsVer.Range("A6").Value = "1.0.4"
$wsVer.Range("C6").Value = "Added the order status widget along the top"
[void]$wsVer.Range("B6").Select()

# --- Restore "Bug List" as the selected/active tab --------------------------
[void]$wsBugs.Activate()
